$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.74884953556969558
$ws.Range("BO1").Value = 0.97895516809551575
$ws.Range("A2").Value = 0.88028813607197476
$ws.Range("C2").Value = 0.81592098434704274
$ws.Range("BP2").Value = 0.92865294479889315
$ws.Range("C4").Value = 0.9103404817095575
$ws.Range("F4").Value = 0.9745257915463732
$ws.Range("C5").Value = 0.7806962025863402
$ws.Range("D5").Value = 0.60814650179708185
$ws.Range("F5").Value = 0.81426175875646578
$ws.Range("AE6").Value = 0.67984731576916724
$ws.Range("E7").Value = 0.97386655391209209
$ws.Range("F7").Value = 0.97968263179143533
$ws.Range("F8").Value = 0.89536582436292711
$ws.Range("J8").Value = 0.91388925893458506
$ws.Range("BJ8").Value = 0.72891380137112505
$ws.Range("G9").Value = 0.98328992103498436
$ws.Range("H9").Value = 0.88230416182286864
$ws.Range("O9").Value = 0.81261199524221461
$ws.Range("I10").Value = 0.99626292173355002
$ws.Range("AV10").Value = 0.88196874389754454
$ws.Range("I11").Value = 0.51141900483763703
$ws.Range("M11").Value = 0.88763589772238838
$ws.Range("AI11").Value = 0.84858821692046105
$ws.Range("J12").Value = 0.97635738810555017
$ws.Range("X12").Value = 0.77022351994630922
$ws.Range("L13").Value = 0.91688942957056596
$ws.Range("L14").Value = 0.76683328880559476
$ws.Range("M14").Value = 0.9897108967548035
$ws.Range("Z14").Value = 0.84440429754074864
$ws.Range("BC14").Value = 0.83374421065590865
$ws.Range("N15").Value = 0.90545645188534885
$ws.Range("B16").Value = 0.68404367345816375
$ws.Range("O16").Value = 0.94952909569747845
$ws.Range("BO16").Value = 0.57810015212316612
$ws.Range("P17").Value = 0.83387374775399525
$ws.Range("R17").Value = 0.81389422321398297
$ws.Range("S18").Value = 0.74905239432349435
$ws.Range("AR18").Value = 0.96129136219367228
$ws.Range("Q19").Value = 0.87980787746126687
$ws.Range("T19").Value = 0.96959405346442207
$ws.Range("J20").Value = 0.88270667559024552
$ws.Range("AN20").Value = 0.62307985922830067
$ws.Range("S21").Value = 0.91859956568357171
$ws.Range("W21").Value = 0.99409214393774525
$ws.Range("AM21").Value = 0.67513183003758026
$ws.Range("U22").Value = 0.60479657155705335
$ws.Range("W22").Value = 0.80916843077620149
$ws.Range("AO22").Value = 0.84744150290945797
$ws.Range("X23").Value = 0.86353952351864927
$ws.Range("O24").Value = 0.79968460924113782
$ws.Range("V24").Value = 0.88685160949016728
$ws.Range("W25").Value = 0.88470955193104028
$ws.Range("Z25").Value = 0.61341880706786722
$ws.Range("AF25").Value = 0.82623607910536712
$ws.Range("AH26").Value = 0.84021049731122288
$ws.Range("Z27").Value = 0.91740789996814931
$ws.Range("Z28").Value = 0.85478929268296855
$ws.Range("AA28").Value = 0.88075598991781501
$ws.Range("BJ28").Value = 0.98160723005804096
$ws.Range("AE29").Value = 0.82763092824071083
$ws.Range("AB30").Value = 0.62089235975782442
$ws.Range("AC30").Value = 0.93506054533390048
$ws.Range("AD31").Value = 0.66018397108831994
$ws.Range("AM31").Value = 0.84004438585807628
$ws.Range("BL32").Value = 0.76576709082997052
$ws.Range("AF34").Value = 0.80547323764031886
$ws.Range("AG34").Value = 0.96768462508815656
$ws.Range("AI34").Value = 0.79805833152070293
$ws.Range("AY34").Value = 0.96140080215296353
$ws.Range("F35").Value = 0.78883276506022426
$ws.Range("M35").Value = 0.88978441422576782
$ws.Range("AG35").Value = 0.98252244681035172
$ws.Range("AK36").Value = 0.95163189083796307
$ws.Range("AJ38").Value = 0.88566899895533635
$ws.Range("AK38").Value = 0.96124573796248947
$ws.Range("AR38").Value = 0.93818845213214308
$ws.Range("BK38").Value = 0.90519423541718103
$ws.Range("AK39").Value = 0.99128081529063228
$ws.Range("AL39").Value = 0.90786158820041529
$ws.Range("AM40").Value = 0.71244581861380141
$ws.Range("AP40").Value = 0.85363009979059989
$ws.Range("AQ41").Value = 0.70733383944996531
$ws.Range("AQ42").Value = 0.80068379886368468
$ws.Range("BN42").Value = 0.72346380161225765
$ws.Range("AR43").Value = 0.9399981600753875
$ws.Range("AP44").Value = 0.82943198232041981
$ws.Range("AQ45").Value = 0.79707848862852304
$ws.Range("AT45").Value = 0.9210753027380183
$ws.Range("AU45").Value = 0.83933430731834002
$ws.Range("AR46").Value = 0.95985092217327628
$ws.Range("AN47").Value = 0.97952830314276318
$ws.Range("AT47").Value = 0.82836333651182747
$ws.Range("AV47").Value = 0.71677458175192932
$ws.Range("AT48").Value = 0.81777676222035955
$ws.Range("AW48").Value = 0.77533974709501952
$ws.Range("G49").Value = 0.88298995077017683
$ws.Range("AD49").Value = 0.92303704732545588
$ws.Range("AU49").Value = 0.62711263916635696
$ws.Range("AY49").Value = 0.90492474523823874
$ws.Range("BB49").Value = 0.94005222434706581
$ws.Range("AV50").Value = 0.84545317584406132
$ws.Range("AZ50").Value = 0.71880284833866748
$ws.Range("B51").Value = 0.97617689643601868
$ws.Range("AX51").Value = 0.77574909301737272
$ws.Range("AA52").Value = 0.95410992333743838
$ws.Range("BK53").Value = 0.95895495279412502
$ws.Range("AK54").Value = 0.95195822992185175
$ws.Range("AZ54").Value = 0.80324224975041236
$ws.Range("BA54").Value = 0.74471702431193587
$ws.Range("M55").Value = 0.76118955266150878
$ws.Range("BA55").Value = 0.6880647017430459
$ws.Range("BE55").Value = 0.93942285348606003
$ws.Range("BC56").Value = 0.76958057921778722
$ws.Range("BE56").Value = 0.73649597311487081
$ws.Range("BF56").Value = 0.78562681655378108
$ws.Range("AD57").Value = 0.60184311641251864
$ws.Range("BE58").Value = 0.93470274252859298
$ws.Range("BG58").Value = 0.9870666610626444
$ws.Range("BH58").Value = 0.80872335569656695
$ws.Range("AH59").Value = 0.91975665956494224
$ws.Range("K60").Value = 0.82885923364421821
$ws.Range("BG60").Value = 0.88197017734233174
$ws.Range("BI60").Value = 0.59832580336960217
$ws.Range("BJ61").Value = 0.94560092740840729
$ws.Range("BK61").Value = 0.95777234830156677
$ws.Range("BH62").Value = 0.835357922372163
$ws.Range("BJ63").Value = 0.70075990177281833
$ws.Range("BL63").Value = 0.99455322839540194
$ws.Range("P64").Value = 0.9957384360399284
$ws.Range("BL65").Value = 0.9680781119434696
$ws.Range("BN65").Value = 0.80538367010878942
$ws.Range("BO65").Value = 0.76079731451249433
$ws.Range("X66").Value = 0.79868923383879442
$ws.Range("BL66").Value = 0.75675590177015795
$ws.Range("AA68").Value = 0.82419368134407245
